# Bugfix: merge the "<Address>" cell (B4) together with the "<SQL Text>"
# cell (previously merged C4:I4) into a single B4:I4 merged cell, giving
# more room to display the SQL statement on the "Execution Plans" sheet.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Item("Execution Plans")

$ws.Range("B4:C4").UnMerge()
$ws.Range("B4:I4").Merge()
$ws.Range("B4").Value = "<SQL Text>"

# Reflect the new selection on the sheet, then restore the first sheet
# (Delta V$SQLAREA) as the active/selected tab, matching the original file.
$ws.Range("B4:I4").Select()
$ws1.Activate()
